$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data refresh: updated coin prices / 1h volume changes, and two row swaps
# (InjectiveProtocol/Toncoin, Celestia/ARBITRUM, Stacks/RocketPoolETH->WOONetwork/Stacks)

$ws.Range("D2").Value = "'43.749.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "'2.296.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'97.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.64%  "
$ws.Range("D6").Value = "'268.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").Value = "'0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("D10").Value = "'45.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("D11").Value = "'0.0934"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "'7.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").Value = "'0.106"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").Value = "'15.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.62%  "
$ws.Range("D15").Value = "'2.640.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "'0.858"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "'2.298.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("D18").Value = "'43.741.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "'0.0000111"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.84%  "
$ws.Range("D20").Value = "'6.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("D21").Value = "'72.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("D22").Value = "'2.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.27%  "
$ws.Range("D23").Value = "'233.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.91%  "
$ws.Range("D24").Value = "'9.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.29%  "
$ws.Range("E25").Value = "  +5.78%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "'11.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("E28").Value = "  +2.09%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'38.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.86%  "
$ws.Range("D31").Value = "'175.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.99%  "
$ws.Range("D32").Value = "'21.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("D33").Value = "'0.0903"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("E34").Value = "  -2.45%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "'4.53"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.08%  "
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("D38").Value = "'0.0352"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("D39").Value = "'3.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.83%  "
$ws.Range("D40").Value = "'0.240"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.64%  "
$ws.Range("D41").Value = "'2.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.63%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'1.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").Value = "'12.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.80%  "
$ws.Range("D44").Value = "'64.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.47%  "
$ws.Range("D45").Value = "'8.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.17%  "
$ws.Range("E46").Value = "  -5.02%  "
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").Value = "'97.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.88%  "
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").Value = "'0.435"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.60%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.48%  "
